# Adds final survival data for resazurin (rows 402-441) to the "data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newRows = @(
    @(20240815,76,"fresh-water-temperature","18C",1,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",2,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",3,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",4,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",5,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",6,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",7,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",8,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",9,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",10,0,0,0,0,1,1),
    @(20240815,76,"fresh-water-temperature","18C",11,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",12,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",13,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",14,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",15,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",16,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",17,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",18,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","18C",19,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","18C",20,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",21,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",22,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",23,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",24,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",25,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",26,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",27,0,0,0,0,0,1),
    @(20240815,76,"fresh-water-temperature","42C",28,0,0,0,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",29,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",30,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",31,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",32,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",33,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",34,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",35,0,0,0,0,1,1),
    @(20240815,76,"fresh-water-temperature","42C",36,0,0,0,0,0,0),
    @(20240815,76,"fresh-water-temperature","42C",37,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",38,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",39,0,0,1,1,1,1),
    @(20240815,76,"fresh-water-temperature","42C",40,0,0,0,0,0,1)
)

$startRow = 402
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# Move the selection to reflect where the author's cursor ended up after
# entering the new data (the frozen header row stays frozen as before).
$ws.Activate() | Out-Null
$ws.Range("L435").Select() | Out-Null
